# Bump benchmark results to v.02: add a third results column (C) with the
# new timings, copying the cell formatting from the matching column B cell
# in each row so the look stays consistent with the existing data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for the new column
$ws.Range("B1").Copy()
$ws.Range("C1").PasteSpecial(-4122)
$ws.Range("C1").Value = "v.02"

# Data rows: copy format from column B (same row) into column C, then set value
$ws.Range("B2").Copy()
$ws.Range("C2").PasteSpecial(-4122)
$ws.Range("C2").Value = "0 ms"

$ws.Range("B3").Copy()
$ws.Range("C3").PasteSpecial(-4122)
$ws.Range("C3").Value = "9878 ms"

$ws.Range("B4").Copy()
$ws.Range("C4").PasteSpecial(-4122)
$ws.Range("C4").Value = "27 ms"

$ws.Range("B5").Copy()
$ws.Range("C5").PasteSpecial(-4122)
$ws.Range("C5").Value = "5294 ms"

$ws.Range("B6").Copy()
$ws.Range("C6").PasteSpecial(-4122)
$ws.Range("C6").Value = "12353 ms"

$ws.Range("B7").Copy()
$ws.Range("C7").PasteSpecial(-4122)
$ws.Range("C7").Value = "75 ms"

# Row 8 (C8) is left untouched - no new timing value for this row.

# Row 9 (WholeTest total) only picks up formatting, no value yet.
$ws.Range("B9").Copy()
$ws.Range("C9").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# Update the selected cell to match the saved view state
$ws.Range("K12").Select()
